$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the 3-bit opcode column (B3:B9) -- introduces new shared strings 000..110
$ws.Range("B3").Value2 = "000"
$ws.Range("B4").Value2 = "001"
$ws.Range("B5").Value2 = "010"
$ws.Range("B6").Value2 = "011"
$ws.Range("B7").Value2 = "100"
$ws.Range("B8").Value2 = "101"
$ws.Range("B9").Value2 = "110"

# Rows that only had a single microinstruction entry move it from column E
# into column F (rows 5, 6, 9 already use both E and F and are untouched).
$ws.Range("F4").Value2 = $ws.Range("E4").Value2
$ws.Range("E4").ClearContents()

$ws.Range("F7").Value2 = $ws.Range("E7").Value2
$ws.Range("E7").ClearContents()

$ws.Range("F8").Value2 = $ws.Range("E8").Value2
$ws.Range("E8").ClearContents()

# Column F widens to match column E (closest width achievable through the
# ColumnWidth property).
$ws.Columns.Item(6).ColumnWidth = 25.5

# Update the active selection to F8
$ws.Range("F8").Select()

# Shrink the workbook window
$excel.ActiveWindow.Width = 7500
$excel.ActiveWindow.Height = 4860
